$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the test-name cells to reflect the new "contrast" tests
$ws.Range("A2").Value = "testDarkContrast"
$ws.Range("A3").Value = "testLightContrast"

# Remove the now-unused rows 4 and 5 (test3/test4 + FAIL entries)
$ws.Rows.Item(4).Delete()
$ws.Rows.Item(4).Delete()

# Re-fit column A's width now that the text is shorter
$ws.Columns.Item(1).AutoFit()
